# Updates crypto price/volume data per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like numbers (e.g. "1.000", "25.511.76").
# Force it to stay Text so Excel does not auto-convert it to a numeric value,
# matching the original inlineStr (text) cells in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.511.76"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "1.666.92"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "234.32"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4654"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").Value = "0.2576"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "0.06138"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "1.666.93"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "0.06971"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "14.62"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "4.340"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "74.90"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.5719"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "25.508.56"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("D19").Value = "0.000006691"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "11.35"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "1.881.36"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "4.420"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").Value = "8.713"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D24").Value = "5.205"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "136.24"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "14.94"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "1.379"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "104.17"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.701"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").Value = "3.973"
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("D31").Value = "0.07774"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "3.602"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "0.04280"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "2.628"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").Value = "0.9452"
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").Value = "0.5976"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").Value = "0.9306"
$ws.Range("E37").Value = "  +15.22%  "
$ws.Range("D38").Value = "2.513"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "1.001"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "101.34"
$ws.Range("E40").Value = "  +4.60%  "
$ws.Range("D41").Value = "0.01468"
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("D42").Value = "1.821"
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("D43").Value = "0.3705"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").Value = "4.912"
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("D45").Value = "0.1105"
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.05258"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "6.123"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("D48").Value = "29.62"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "7.377"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "0.9993"
$ws.Range("E51").Value = "  +0.33%  "
